# Ajout de tache fonctionnel + Fonctionnement Theorique sous linux
# Appends the new "taches" (tasks) entered by the user to column A,
# rows 8 through 27, on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$taches = @(
    "Test",
    "Saisie",
    "test5",
    "TEST55",
    "Test77",
    "1234",
    "Test112",
    "tet",
    "tet4",
    "test99",
    "test100",
    "test101",
    "123456",
    "voici ma tâche",
    "voici ma tâche",
    "voici ma tâche",
    "voici ma tâche",
    "voici ma tâche2",
    "voici ma tâche 3",
    "test55555"
)

$startRow = 8
for ($i = 0; $i -lt $taches.Count; $i++) {
    $row = $startRow + $i
    $value = $taches[$i]
    $cell = $ws.Cells.Item($row, 1)

    # Some entries look like plain numbers ("1234", "123456"); they were
    # typed as free-text task names, so force the cell to text formatting
    # before writing the value to keep it from being stored as a number.
    if ($value -match '^[0-9]+$') {
        $cell.NumberFormat = "@"
    }

    $cell.Value = $value
}
